# The workbook holds a weekly "Apio" (celery) price-reporting log where each
# row is one market observation, ordered by date. This commit adds one more
# weekly observation. In the canonical OOXML this shows up as every existing
# row from 28 downward being pushed one row lower (28->29, 29->30, ... ,
# 122->123) and a brand-new row being written into the now-vacant row 28.
#
# Reproduce that with a native row insert at row 28 (which shifts everything
# below it down by one, carrying formatting along) and then fill the new
# row 28 with the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 28; rows 28..122 become 29..123.
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new weekly record.
$ws.Cells.Item(28, 1).Value = 5
$ws.Cells.Item(28, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(28, 3).Value = "Maule"
$ws.Cells.Item(28, 4).Value = 44481
$ws.Cells.Item(28, 5).Value = 7
$ws.Cells.Item(28, 6).Value = 100112017
$ws.Cells.Item(28, 7).Value = "Apio"
$ws.Cells.Item(28, 8).Value = "Americana (o)"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 600
$ws.Cells.Item(28, 11).Value = 7500
$ws.Cells.Item(28, 12).Value = 7500
$ws.Cells.Item(28, 13).Value = 7500
$ws.Cells.Item(28, 14).Value = "`$/docena de matas"
$ws.Cells.Item(28, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(28, 16).Value = 1250
$ws.Cells.Item(28, 17).Value = 6
$ws.Cells.Item(28, 18).Value = "Hortaliza"
